$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 81
$ws.Range("F6").Value = 868
$ws.Range("F7").Value = 448
$ws.Range("F8").Value = 4763
$ws.Range("F9").Value = 4763
$ws.Range("F10").Value = 109
$ws.Range("F11").Value = 124
$ws.Range("F15").Value = 126
$ws.Range("F16").Value = 7624
$ws.Range("F21").Value = 535
$ws.Range("F24").Value = 6287
$ws.Range("F25").Value = 2257
$ws.Range("F26").Value = 22
$ws.Range("F27").Value = 2093
$ws.Range("F28").Value = 6204
$ws.Range("F30").Value = 31
$ws.Range("F31").Value = 119
$ws.Range("F32").Value = 91
$ws.Range("F34").Value = 6502
$ws.Range("F36").Value = 210
$ws.Range("F39").Value = 21
$ws.Range("F43").Value = 62
$ws.Range("F44").Value = 1117
$ws.Range("F47").Value = 2155

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 235
$ws.Range("F6").Value = 130

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 235
$ws.Range("F7").Value = 81
$ws.Range("F9").Value = 448
$ws.Range("F10").Value = 4763
$ws.Range("F11").Value = 4763
$ws.Range("F12").Value = 109
$ws.Range("F13").Value = 124
$ws.Range("F17").Value = 126
$ws.Range("F18").Value = 7624
$ws.Range("F19").Value = 7624
$ws.Range("F22").Value = 535
$ws.Range("F24").Value = 130
$ws.Range("F25").Value = 6287
$ws.Range("F26").Value = 2257
$ws.Range("F27").Value = 2093
$ws.Range("F29").Value = 6204
$ws.Range("F32").Value = 31
$ws.Range("F33").Value = 119
$ws.Range("F34").Value = 91
$ws.Range("F36").Value = 6503
$ws.Range("F38").Value = 210
$ws.Range("F41").Value = 21
$ws.Range("F45").Value = 1117
$ws.Range("F49").Value = 2155

